# Updated cryptos list (Price / Volume(1h) columns) per GitHub Actions refresh.
# Numeric-looking Price values are prefixed with a leading apostrophe so Excel
# stores them as text (matching the original inline-string "Price" column)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.942.16"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.660.24"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("D5").Value = "'217.39"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'28.89"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'0.0903"
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").Value = "1.898.89"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").Value = "1.669.22"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("D14").Value = "'0.610"
$ws.Range("E14").Value = "  +6.46%  "
$ws.Range("D15").Value = "'10.15"
$ws.Range("E15").Value = "  +14.06%  "
$ws.Range("D17").Value = "29.970.72"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "'65.05"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "'243.19"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("E22").Value = "  +4.20%  "
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "'158.80"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").Value = "'15.81"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("D33").Value = "'3.23"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").Value = "1.447.14"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("E35").Value = "  +4.58%  "
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "'79.68"
$ws.Range("E37").Value = "  +14.28%  "
$ws.Range("E38").Value = "  +2.73%  "
$ws.Range("D39").Value = "'0.576"
$ws.Range("E39").Value = "  +3.31%  "
$ws.Range("E41").Value = "  -8.04%  "
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("D43").Value = "'0.0499"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").Value = "'50.57"
$ws.Range("E47").Value = "  -7.12%  "
$ws.Range("D48").Value = "1.804.67"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").Value = "'94.32"
$ws.Range("E50").Value = "  +6.18%  "
$ws.Range("E51").Value = "  +3.50%  "
